$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new customer ("LITTLE BIRD DELICATESSEN") was added to the leaderboard.
# It belongs right after "TAQUERIA Y MERCADO ANDALE" (row 39) and before
# "VALLEY OFFICE PARK" (old row 40), so insert a new row there, pushing
# everything from the old row 40 onward down by one.
$ws.Rows(40).Insert()

# The newly-added customer now has a recent invoice, so the previous row
# (TAQUERIA Y MERCADO ANDALE) picks up a "Last Invoice Date" -- copy the
# date number-format/alignment from a neighboring date cell so the new
# value renders the same way, then set the value.
$ws.Range("D41").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("D39").Value = 45916

# Populate the new row with the new customer's data.
$ws.Range("A40").Value = "LITTLE BIRD DELICATESSEN"
$ws.Range("B40").Value = "House Account"
$ws.Range("C40").Value = "030"
$ws.Range("E40").Value = "0008303"
